$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 176, pushing existing rows 176-183 down to 177-184
$ws.Rows.Item(176).Insert()

# Populate the new row 176 with data (same boilerplate columns A,B,C,E-K as row 177 now contains,
# but new values for D, L, M, N, O, P, Q, R, S, T)
$ws.Cells.Item(176, 1).Value = 3
$ws.Cells.Item(176, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(176, 3).Value = "Coquimbo"
$ws.Cells.Item(176, 4).Value = 44585
$ws.Cells.Item(176, 5).Value = 5
$ws.Cells.Item(176, 6).Value = "Fruta"
$ws.Cells.Item(176, 7).Value = 100101
$ws.Cells.Item(176, 8).Value = "Berries"
$ws.Cells.Item(176, 9).Value = 100101001
$ws.Cells.Item(176, 10).Value = "Arándano (blue)"
$ws.Cells.Item(176, 11).Value = "Sin especificar"
$ws.Cells.Item(176, 12).Value = "Primera"
$ws.Cells.Item(176, 13).Value = 220
$ws.Cells.Item(176, 14).Value = 4000
$ws.Cells.Item(176, 15).Value = 4500
$ws.Cells.Item(176, 16).Value = 4273
$ws.Cells.Item(176, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(176, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(176, 19).Value = 2136
$ws.Cells.Item(176, 20).Value = 2
